$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot each data row (2-31) as a Formula array BEFORE any writes. ---
# Using Formula (not Value) preserves the literal HYPERLINK(...) formula text for columns
# S/T/U/V/W/X/Y/Z while still capturing literal numbers/strings for the rest of the row.
# A snapshot is required because the re-ordering below is a permutation with cycles, so naive
# row-by-row copying would overwrite source rows before they are read.
$snapshot = @{}
$snapshot[2] = $ws.Rows.Item(2).Formula
$snapshot[3] = $ws.Rows.Item(3).Formula
$snapshot[4] = $ws.Rows.Item(4).Formula
$snapshot[5] = $ws.Rows.Item(5).Formula
$snapshot[6] = $ws.Rows.Item(6).Formula
$snapshot[7] = $ws.Rows.Item(7).Formula
$snapshot[8] = $ws.Rows.Item(8).Formula
$snapshot[9] = $ws.Rows.Item(9).Formula
$snapshot[10] = $ws.Rows.Item(10).Formula
$snapshot[11] = $ws.Rows.Item(11).Formula
$snapshot[12] = $ws.Rows.Item(12).Formula
$snapshot[13] = $ws.Rows.Item(13).Formula
$snapshot[14] = $ws.Rows.Item(14).Formula
$snapshot[15] = $ws.Rows.Item(15).Formula
$snapshot[16] = $ws.Rows.Item(16).Formula
$snapshot[17] = $ws.Rows.Item(17).Formula
$snapshot[18] = $ws.Rows.Item(18).Formula
$snapshot[19] = $ws.Rows.Item(19).Formula
$snapshot[20] = $ws.Rows.Item(20).Formula
$snapshot[21] = $ws.Rows.Item(21).Formula
$snapshot[22] = $ws.Rows.Item(22).Formula
$snapshot[23] = $ws.Rows.Item(23).Formula
$snapshot[24] = $ws.Rows.Item(24).Formula
$snapshot[25] = $ws.Rows.Item(25).Formula
$snapshot[26] = $ws.Rows.Item(26).Formula
$snapshot[27] = $ws.Rows.Item(27).Formula
$snapshot[28] = $ws.Rows.Item(28).Formula
$snapshot[29] = $ws.Rows.Item(29).Formula
$snapshot[30] = $ws.Rows.Item(30).Formula
$snapshot[31] = $ws.Rows.Item(31).Formula

# --- Step 2: write every row back out in its new position. ---
# (mapping: new row N <- old row M, derived from the upstream data refresh)
$ws.Rows.Item(2).Formula = $snapshot[2]
$ws.Rows.Item(3).Formula = $snapshot[3]
$ws.Rows.Item(4).Formula = $snapshot[6]  # moved
$ws.Rows.Item(5).Formula = $snapshot[5]
$ws.Rows.Item(6).Formula = $snapshot[7]  # moved
$ws.Rows.Item(7).Formula = $snapshot[4]  # moved
$ws.Rows.Item(8).Formula = $snapshot[8]
$ws.Rows.Item(9).Formula = $snapshot[9]
$ws.Rows.Item(10).Formula = $snapshot[10]
$ws.Rows.Item(11).Formula = $snapshot[11]
$ws.Rows.Item(12).Formula = $snapshot[15]  # moved
$ws.Rows.Item(13).Formula = $snapshot[13]
$ws.Rows.Item(14).Formula = $snapshot[30]  # moved
$ws.Rows.Item(15).Formula = $snapshot[29]  # moved
$ws.Rows.Item(16).Formula = $snapshot[22]  # moved
$ws.Rows.Item(17).Formula = $snapshot[23]  # moved
$ws.Rows.Item(18).Formula = $snapshot[27]  # moved
$ws.Rows.Item(19).Formula = $snapshot[25]  # moved
$ws.Rows.Item(20).Formula = $snapshot[14]  # moved
$ws.Rows.Item(21).Formula = $snapshot[24]  # moved
$ws.Rows.Item(22).Formula = $snapshot[26]  # moved
$ws.Rows.Item(23).Formula = $snapshot[16]  # moved
$ws.Rows.Item(24).Formula = $snapshot[20]  # moved
$ws.Rows.Item(25).Formula = $snapshot[12]  # moved
$ws.Rows.Item(26).Formula = $snapshot[17]  # moved
$ws.Rows.Item(27).Formula = $snapshot[31]  # moved
$ws.Rows.Item(28).Formula = $snapshot[28]
$ws.Rows.Item(29).Formula = $snapshot[19]  # moved
$ws.Rows.Item(30).Formula = $snapshot[18]  # moved
$ws.Rows.Item(31).Formula = $snapshot[21]  # moved

# --- Step 3: the "Förändrad" (changed) column C advances by one day for every data row. ---
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 3).Value = 46077
}

# --- Step 4: undo the engine's autofit-on-write side effect. ---
# Re-writing a whole row's Formula array (Step 2) makes rows containing wrapped text (column R)
# auto-grow to fit their new content, and marks every touched row as having a "custom" height even
# when it never had one. Neither effect is part of the source edit, so put row heights back:
#  - rows 2-30 originally had an explicit 15pt custom height -> restore it explicitly.
#  - row 31 originally had no explicit height at all -> AutoFit() clears the custom-height flag.
for ($r = 2; $r -le 30; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}
$ws.Rows.Item(31).EntireRow.AutoFit()

